$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 187, shifting existing rows 187:216 down to 188:217
$ws.Rows.Item(187).Insert()

# Fill in the new row 187 with the inserted record (weekly price data point)
$ws.Cells.Item(187, 1).Value = 10
$ws.Cells.Item(187, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(187, 3).Value = "La Araucanía"
$ws.Cells.Item(187, 4).Value = 45127
$ws.Cells.Item(187, 5).Value = 9
$ws.Cells.Item(187, 6).Value = 100112031
$ws.Cells.Item(187, 7).Value = "Poroto verde"
$ws.Cells.Item(187, 8).Value = "Sin especificar"
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 60
$ws.Cells.Item(187, 11).Value = 24000
$ws.Cells.Item(187, 12).Value = 24000
$ws.Cells.Item(187, 13).Value = 24000
$ws.Cells.Item(187, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(187, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(187, 16).Value = 960
$ws.Cells.Item(187, 17).Value = 25
$ws.Cells.Item(187, 18).Value = "Hortaliza"
